$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates (shared-string rich text) ---
$ws.Range("A8").Value = "Volume 31   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/2/2024  Through  9/8/2024"

# --- Cells changing from placeholder text ("0") to numeric values ---
$ws.Range("C16").NumberFormat = $ws.Range("D16").NumberFormat
$ws.Range("C16").Value = 1
$ws.Range("C22").NumberFormat = $ws.Range("J22").NumberFormat
$ws.Range("C22").Value = 1
$ws.Range("F22").NumberFormat = $ws.Range("I22").NumberFormat
$ws.Range("F22").Value = 1

# --- Cells changing from numeric values to placeholder text ("0") ---
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0"

# --- Plain numeric value updates ---
$ws.Range("N14").Value = -50
$ws.Range("J15").Value = 14
$ws.Range("K15").Value = -14.285714285714
$ws.Range("L15").Value = -14.285714285714
$ws.Range("M15").Value = 71.428571428571
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 40
$ws.Range("I16").Value = 78
$ws.Range("J16").Value = 58
$ws.Range("K16").Value = 34.482758620689
$ws.Range("L16").Value = 27.868852459016
$ws.Range("M16").Value = -35.537190082644
$ws.Range("N16").Value = -84.114052953156
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = 14.285714285714
$ws.Range("I17").Value = 147
$ws.Range("J17").Value = 147
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = -2.649006622516
$ws.Range("M17").Value = 33.636363636363
$ws.Range("N17").Value = -35.526315789473
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -7.692307692307
$ws.Range("I18").Value = 82
$ws.Range("J18").Value = 84
$ws.Range("K18").Value = -2.380952380952
$ws.Range("L18").Value = -36.923076923076
$ws.Range("M18").Value = -68.461538461538
$ws.Range("N18").Value = -93.027210884353
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 6.666666666666
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -11.538461538461
$ws.Range("I19").Value = 372
$ws.Range("J19").Value = 395
$ws.Range("K19").Value = -5.822784810126
$ws.Range("L19").Value = -13.084112149532
$ws.Range("M19").Value = 36.263736263736
$ws.Range("N19").Value = -16.591928251121
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -15.384615384615
$ws.Range("J20").Value = 105
$ws.Range("K20").Value = 17.142857142857
$ws.Range("L20").Value = 59.740259740259
$ws.Range("N20").Value = -90.096618357487
$ws.Range("C21").Value = 25
$ws.Range("E21").Value = -3.846153846153
$ws.Range("F21").Value = 92
$ws.Range("G21").Value = 100
$ws.Range("H21").Value = -8
$ws.Range("I21").Value = 818
$ws.Range("J21").Value = 805
$ws.Range("K21").Value = 1.614906832298
$ws.Range("L21").Value = -5.104408352668
$ws.Range("M21").Value = -5.868814729574
$ws.Range("N21").Value = -77.309292649098
$ws.Range("I22").Value = 10
$ws.Range("K22").Value = 25
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -37.5
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 14.285714285714
$ws.Range("F24").Value = 70
$ws.Range("G24").Value = 79
$ws.Range("H24").Value = -11.392405063291
$ws.Range("I24").Value = 674
$ws.Range("J24").Value = 760
$ws.Range("K24").Value = -11.315789473684
$ws.Range("L24").Value = -13.256113256113
$ws.Range("M24").Value = 7.667731629392
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = -25.925925925925
$ws.Range("I25").Value = 168
$ws.Range("J25").Value = 227
$ws.Range("K25").Value = -25.991189427312
$ws.Range("L25").Value = -26.95652173913
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 71.428571428571
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = 45.833333333333
$ws.Range("I26").Value = 292
$ws.Range("J26").Value = 239
$ws.Range("K26").Value = 22.175732217573
$ws.Range("L26").Value = 25.862068965517
$ws.Range("M26").Value = -0.341296928327
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = 37.5
$ws.Range("L27").Value = 29.411764705882
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 16.666666666666
$ws.Range("I28").Value = 47
$ws.Range("J28").Value = 49
$ws.Range("K28").Value = -4.081632653061
$ws.Range("L28").Value = -12.962962962963
